$wb = $excel.ActiveWorkbook

# --- zh-cn sheet (row 7 / a2114e04-... handback just landed) ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = "a2114e04-9715-4b53-a924-339af4b35bc3.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/546345a5838a528498261d1923535a9c3d551e21/e2e/a2114e04-9715-4b53-a924-339af4b35bc3.md", [Type]::Missing, [Type]::Missing, "a2114e04-9715-4b53-a924-339af4b35bc3.md") | Out-Null

$wsZh.Range("J7").Value = "a2114e04-9715-4b53-a924-339af4b35bc3.3a97fdb82089c30aca977531c73bd2012b718cbe.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-17 06:48:53"
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98a0976c1aff4d7b4ce5dc501343615d1206f02a/e2e/a2114e04-9715-4b53-a924-339af4b35bc3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/546345a5838a528498261d1923535a9c3d551e21/e2e/a2114e04-9715-4b53-a924-339af4b35bc3.md."

# --- de-de sheet (same handback event) ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = "a2114e04-9715-4b53-a924-339af4b35bc3.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/546345a5838a528498261d1923535a9c3d551e21/e2e/a2114e04-9715-4b53-a924-339af4b35bc3.md", [Type]::Missing, [Type]::Missing, "a2114e04-9715-4b53-a924-339af4b35bc3.md") | Out-Null

$wsDe.Range("J7").Value = "a2114e04-9715-4b53-a924-339af4b35bc3.3a97fdb82089c30aca977531c73bd2012b718cbe.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-17 06:49:00"
$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/98a0976c1aff4d7b4ce5dc501343615d1206f02a/e2e/a2114e04-9715-4b53-a924-339af4b35bc3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/546345a5838a528498261d1923535a9c3d551e21/e2e/a2114e04-9715-4b53-a924-339af4b35bc3.md."
